$d = $word.ActiveDocument

# The document starts with 23 paragraphs:
#   1-19  UserInput ... dialogeCharacter misspelling of dialogue   (removed)
#   20    IntroScreen                                              (kept)
#   21    currentProfile is Any type and should be specified       (kept)
#   22-23 ChatScreen / Chats missing camelCase                     (removed)
#
# Remove everything except the "IntroScreen" heading and the
# "currentProfile ... is Any type and should be specified" paragraph.

# Delete the trailing paragraphs (ChatScreen, Chats missing camelCase) first
# so the earlier paragraph indices/offsets used below stay valid.
$startTail = $d.Paragraphs.Item(21).Range.End
$endTail = $d.Paragraphs.Item($d.Paragraphs.Count).Range.End
$d.Range($startTail, $endTail).Delete()

# Delete the leading paragraphs (UserInput ... dialogeCharacter misspelling of dialogue).
$startHead = $d.Paragraphs.Item(1).Range.Start
$endHead = $d.Paragraphs.Item(19).Range.End
$d.Range($startHead, $endHead).Delete()

# Only two paragraphs remain now: "IntroScreen" (Kop2/Heading 2) and
# "currentProfile is Any type and should be specified" (Geenafstand/No Spacing).
# The "_GoBack" bookmark currently still sits at the end of the "IntroScreen"
# text (its original position); move it to the end of the second paragraph's
# text instead, mirroring the last edit having happened there.
$p2 = $d.Paragraphs.Item(2)
$textEnd = $p2.Range.End - 1

# Bookmarks.Add with a zero-length range placed exactly at "end of paragraph
# text" (i.e. right before the paragraph mark) snaps back to paragraph 1
# instead of honouring the requested location. Work around it by briefly
# inserting a placeholder character to shift that boundary, adding the
# bookmark at the now safe mid-text offset, then deleting the placeholder.
$d.Range($textEnd, $textEnd).InsertAfter("~")
$d.Bookmarks.Add("_GoBack", $d.Range($textEnd, $textEnd))
$d.Range($textEnd, $textEnd + 1).Delete()
